# Update cryptos list (prices, volume %, and row 47/48 coin swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.891.48"
$ws.Range("E2").Value = "  -3.38%  "

# Row 3
$ws.Range("D3").Value = "2.488.29"
$ws.Range("E3").Value = "  -5.96%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "558.15"
$ws.Range("E5").Value = "  -4.05%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.52"
$ws.Range("E6").Value = "  -5.15%  "

# Row 7
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.604"
$ws.Range("E8").Value = "  -3.69%  "

# Row 9
$ws.Range("D9").Value = "2.485.36"
$ws.Range("E9").Value = "  -6.01%  "

# Row 10
$ws.Range("E10").Value = "  -8.01%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.50"
$ws.Range("E11").Value = "  -5.65%  "

# Row 12
$ws.Range("E12").Value = "  -1.41%  "

# Row 13
$ws.Range("E13").Value = "  -6.24%  "

# Row 14
$ws.Range("E14").Value = "  -6.81%  "

# Row 15
$ws.Range("D15").Value = "2.936.96"
$ws.Range("E15").Value = "  -6.03%  "

# Row 16
$ws.Range("E16").Value = "  -8.23%  "

# Row 17
$ws.Range("D17").Value = "61.776.72"
$ws.Range("E17").Value = "  -3.47%  "

# Row 18
$ws.Range("D18").Value = "2.480.80"
$ws.Range("E18").Value = "  -6.23%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.30"
$ws.Range("E19").Value = "  -7.65%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.23"
$ws.Range("E20").Value = "  -6.85%  "

# Row 21
$ws.Range("E21").Value = "  -6.19%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "323.84"
$ws.Range("E22").Value = "  -6.66%  "

# Row 23
$ws.Range("E23").Value = "  -0.01%  "

# Row 24
$ws.Range("E24").Value = "  +3.04%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.68"
$ws.Range("E25").Value = "  -5.44%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000102"
$ws.Range("E26").Value = "  -9.39%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "575.99"
$ws.Range("E27").Value = "  -2.17%  "

# Row 28
$ws.Range("D28").Value = "2.604.78"
$ws.Range("E28").Value = "  -6.26%  "

# Row 29
$ws.Range("E29").Value = "  -6.41%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.00%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.43"
$ws.Range("E31").Value = "  -10.17%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.85"
$ws.Range("E32").Value = "  -4.67%  "

# Row 33
$ws.Range("E33").Value = "  -6.46%  "

# Row 34
$ws.Range("E34").Value = "  -6.46%  "

# Row 35
$ws.Range("E35").Value = "  -7.97%  "

# Row 36
$ws.Range("E36").Value = "  -9.44%  "

# Row 37
$ws.Range("E37").Value = "  -9.19%  "

# Row 38
$ws.Range("E38").Value = "  -0.09%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.387"
$ws.Range("E39").Value = "  -4.50%  "

# Row 40
$ws.Range("E40").Value = "  -5.56%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "145.02"
$ws.Range("E41").Value = "  -4.43%  "

# Row 42
$ws.Range("E42").Value = "  -6.86%  "

# Row 43
$ws.Range("E43").Value = "  +0.06%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.49"
$ws.Range("E44").Value = "  -4.06%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "40.69"
$ws.Range("E45").Value = "  -2.91%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "149.97"
$ws.Range("E46").Value = "  -8.56%  "

# Row 47
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.22"
$ws.Range("E47").Value = "  -9.59%  "

# Row 48
$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.68"
$ws.Range("E48").Value = "  -6.37%  "

# Row 49
$ws.Range("E49").Value = "  -7.84%  "

# Row 50
$ws.Range("E50").Value = "  -5.50%  "

# Row 51
$ws.Range("E51").Value = "  -5.46%  "
